$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: add the new commit entry text and hours
$ws.Range("C34").Value = "potions working & object classes update"
$ws.Range("G34").Value = 1.8

# Update the total formula to include row 34
$ws.Range("G39").Formula = "=SUM(G4:G34)"

# Update selection to C42
$ws.Range("C42").Select()
